$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "PALACIOS PANTA LUIS MIGUEL",
    "PANTA VARONA CANDY ELIZABETH",
    "MAZA RIOFRIO CINTHIA NATELAHI",
    "SALAZAR VEGA MARIA FERNANDA",
    "PANTA NIMA FREDDY ROLAND JUNIOR",
    "CRISANTO CARMEN ROSITA ABIGAIL",
    "ELIAS MACHADO JUANA MARGOT",
    "HIDALGO MOSCOL YESSICA JAZMIN",
    "VEGA ZAPATA JESUS GABRIEL",
    "TALLEDO ELIAS ANDREA ALESSANDRA"
)

$totals = @(82, 80, 79, 78, 77, 77, 76, 73, 72, 66)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $totals[$i]
}
